$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 7449
$ws.Cells.Item(28, 9).Value = 787.8946999999999
$ws.Cells.Item(28, 10).Value = 32761.2
$ws.Cells.Item(28, 11).Value = 787.8946999999999
$ws.Cells.Item(28, 12).Value = 32761.2
$ws.Cells.Item(28, 13).Value = -302.8946999999999
$ws.Cells.Item(28, 14).Value = -33731.2
$ws.Cells.Item(33, 8).Value = 748.6667
$ws.Cells.Item(33, 9).Value = 838.2
$ws.Cells.Item(33, 11).Value = 838.2
$ws.Cells.Item(33, 13).Value = -609.2
$ws.Cells.Item(94, 8).Value = 1250
$ws.Cells.Item(112, 8).Value = 3890.5
$ws.Cells.Item(112, 10).Value = 3994.4465
$ws.Cells.Item(112, 12).Value = 11983.3395
$ws.Cells.Item(112, 14).Value = -14199.3395
$ws.Cells.Item(129, 8).Value = 1183.7887
$ws.Cells.Item(129, 9).Value = 741.7143
$ws.Cells.Item(129, 10).Value = 1232.1406
$ws.Cells.Item(129, 11).Value = 2225.1429
$ws.Cells.Item(129, 12).Value = 3696.4218
$ws.Cells.Item(129, 13).Value = 2774.8571
$ws.Cells.Item(129, 14).Value = -13696.4218
$ws.Cells.Item(138, 8).Value = 3803.9177
$ws.Cells.Item(138, 10).Value = 3907.6616
$ws.Cells.Item(138, 12).Value = 11722.9848
$ws.Cells.Item(138, 14).Value = -22002.9848
$ws.Cells.Item(140, 8).Value = 72992.94
$ws.Cells.Item(140, 10).Value = 72992.94
$ws.Cells.Item(140, 12).Value = 72992.94
$ws.Cells.Item(140, 14).Value = -83352.94

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6503839
$ws.Cells.Item(32, 9).Value = 7254994
$ws.Cells.Item(32, 10).Value = 25125
$ws.Cells.Item(32, 11).Value = 7254994
$ws.Cells.Item(32, 12).Value = 25125
$ws.Cells.Item(32, 13).Value = -7254707
$ws.Cells.Item(32, 14).Value = -25699
$ws.Cells.Item(34, 8).Value = 100028
$ws.Cells.Item(34, 10).Value = 100028
$ws.Cells.Item(34, 12).Value = 100028
$ws.Cells.Item(34, 14).Value = -100570
$ws.Cells.Item(122, 8).Value = 64201.812
$ws.Cells.Item(122, 9).Value = 78356.08
$ws.Cells.Item(122, 10).Value = 2866.6667
$ws.Cells.Item(122, 11).Value = 235068.24
$ws.Cells.Item(122, 12).Value = 8600.000100000001
$ws.Cells.Item(122, 13).Value = -232618.24
$ws.Cells.Item(122, 14).Value = -13500.0001
$ws.Cells.Item(132, 8).Value = 963697.4399999999
$ws.Cells.Item(132, 9).Value = 1593.0476
$ws.Cells.Item(132, 11).Value = 4779.142800000001
$ws.Cells.Item(132, 13).Value = -2249.142800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 60062.293
$ws.Cells.Item(107, 9).Value = 100901.9
$ws.Cells.Item(107, 10).Value = 1720
$ws.Cells.Item(107, 11).Value = 100901.9
$ws.Cells.Item(107, 12).Value = 1720
$ws.Cells.Item(107, 13).Value = -98981.89999999999
$ws.Cells.Item(107, 14).Value = -5560

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2002.0625
$ws.Cells.Item(16, 9).Value = 1918.8889
$ws.Cells.Item(16, 10).Value = 2109
$ws.Cells.Item(16, 11).Value = 1918.8889
$ws.Cells.Item(16, 12).Value = 2109
$ws.Cells.Item(16, 13).Value = -1631.8889
$ws.Cells.Item(16, 14).Value = -2683
$ws.Cells.Item(28, 8).Value = 70000
$ws.Cells.Item(28, 10).Value = 70000
$ws.Cells.Item(28, 12).Value = 70000
$ws.Cells.Item(28, 14).Value = -70490
$ws.Cells.Item(31, 8).Value = 5910.7915
$ws.Cells.Item(31, 9).Value = 1782.2667
$ws.Cells.Item(31, 10).Value = 6997.2456
$ws.Cells.Item(31, 11).Value = 1782.2667
$ws.Cells.Item(31, 12).Value = 6997.2456
$ws.Cells.Item(31, 13).Value = -1487.2667
$ws.Cells.Item(31, 14).Value = -7587.2456
$ws.Cells.Item(34, 8).Value = 5910.7915
$ws.Cells.Item(34, 9).Value = 1782.2667
$ws.Cells.Item(34, 10).Value = 6997.2456
$ws.Cells.Item(34, 11).Value = 1782.2667
$ws.Cells.Item(34, 12).Value = 6997.2456
$ws.Cells.Item(34, 13).Value = -1580.2667
$ws.Cells.Item(34, 14).Value = -7401.2456
$ws.Cells.Item(113, 8).Value = 2002.0625
$ws.Cells.Item(113, 9).Value = 1918.8889
$ws.Cells.Item(113, 10).Value = 2109
$ws.Cells.Item(113, 11).Value = 1918.8889
$ws.Cells.Item(113, 12).Value = 2109
$ws.Cells.Item(113, 13).Value = 251.1111000000001
$ws.Cells.Item(113, 14).Value = -6449
$ws.Cells.Item(132, 8).Value = 27780354
$ws.Cells.Item(132, 9).Value = 29414138
$ws.Cells.Item(132, 10).Value = 23812596
$ws.Cells.Item(132, 11).Value = 88242414
$ws.Cells.Item(132, 12).Value = 71437788
$ws.Cells.Item(132, 13).Value = -88239884
$ws.Cells.Item(132, 14).Value = -71442848

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 979.5357
$ws.Cells.Item(5, 9).Value = 678.3889
$ws.Cells.Item(5, 10).Value = 1521.6
$ws.Cells.Item(5, 11).Value = 2035.1667
$ws.Cells.Item(5, 12).Value = 4564.799999999999
$ws.Cells.Item(5, 13).Value = -1923.1667
$ws.Cells.Item(5, 14).Value = -4788.799999999999
$ws.Cells.Item(129, 8).Value = 798661.5
$ws.Cells.Item(129, 10).Value = 1213740.6
$ws.Cells.Item(129, 12).Value = 3641221.8
$ws.Cells.Item(129, 14).Value = -3651221.8
$ws.Cells.Item(135, 8).Value = 979.5357
$ws.Cells.Item(135, 9).Value = 678.3889
$ws.Cells.Item(135, 10).Value = 1521.6
$ws.Cells.Item(135, 11).Value = 6105.5001
$ws.Cells.Item(135, 12).Value = 13694.4
$ws.Cells.Item(135, 13).Value = -3570.5001
$ws.Cells.Item(135, 14).Value = -18764.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 2793.1333
$ws.Cells.Item(97, 9).Value = 2838.6667
$ws.Cells.Item(97, 10).Value = 2724.8333
$ws.Cells.Item(97, 11).Value = 2838.6667
$ws.Cells.Item(97, 12).Value = 2724.8333
$ws.Cells.Item(97, 13).Value = -2342.6667
$ws.Cells.Item(97, 14).Value = -3716.8333
$ws.Cells.Item(107, 8).Value = 455.26666
$ws.Cells.Item(107, 9).Value = 318
$ws.Cells.Item(107, 10).Value = 729.8
$ws.Cells.Item(107, 11).Value = 318
$ws.Cells.Item(107, 12).Value = 729.8
$ws.Cells.Item(107, 13).Value = 1602
$ws.Cells.Item(107, 14).Value = -4569.8
$ws.Cells.Item(126, 8).Value = 5245
$ws.Cells.Item(126, 10).Value = 7000
$ws.Cells.Item(126, 12).Value = 21000
$ws.Cells.Item(126, 14).Value = -25940

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 8844.444
$ws.Cells.Item(2, 10).Value = 9324.875
$ws.Cells.Item(2, 12).Value = 9324.875
$ws.Cells.Item(2, 14).Value = -9548.875
$ws.Cells.Item(7, 8).Value = 9741.286
$ws.Cells.Item(7, 9).Value = 9246
$ws.Cells.Item(7, 11).Value = 9246
$ws.Cells.Item(7, 13).Value = -9134
$ws.Cells.Item(16, 8).Value = 566
$ws.Cells.Item(16, 9).Value = 279.91666
$ws.Cells.Item(16, 10).Value = 3999
$ws.Cells.Item(16, 11).Value = 279.91666
$ws.Cells.Item(16, 12).Value = 3999
$ws.Cells.Item(16, 13).Value = -109.91666
$ws.Cells.Item(16, 14).Value = -4339
$ws.Cells.Item(40, 8).Value = 3107.9375
$ws.Cells.Item(40, 9).Value = 2670.923
$ws.Cells.Item(40, 10).Value = 5001.6665
$ws.Cells.Item(40, 11).Value = 2670.923
$ws.Cells.Item(40, 12).Value = 5001.6665
$ws.Cells.Item(40, 13).Value = -2534.923
$ws.Cells.Item(40, 14).Value = -5273.6665
$ws.Cells.Item(46, 10).Value = 500
$ws.Cells.Item(46, 12).Value = 500
$ws.Cells.Item(46, 14).Value = -876
$ws.Cells.Item(68, 8).Value = 2000
$ws.Cells.Item(68, 9).Value = 2000
$ws.Cells.Item(68, 10).Value = 2000
$ws.Cells.Item(68, 11).Value = 2000
$ws.Cells.Item(68, 12).Value = 2000
$ws.Cells.Item(68, 13).Value = -1251
$ws.Cells.Item(68, 14).Value = -3498
$ws.Cells.Item(71, 8).Value = 2000
$ws.Cells.Item(71, 9).Value = 2000
$ws.Cells.Item(71, 10).Value = 2000
$ws.Cells.Item(71, 11).Value = 10000
$ws.Cells.Item(71, 12).Value = 10000
$ws.Cells.Item(71, 13).Value = -6256
$ws.Cells.Item(71, 14).Value = -17488
$ws.Cells.Item(82, 8).Value = 38464204
$ws.Cells.Item(82, 9).Value = 71431290
$ws.Cells.Item(82, 10).Value = 2605.3333
$ws.Cells.Item(82, 11).Value = 71431290
$ws.Cells.Item(82, 12).Value = 2605.3333
$ws.Cells.Item(82, 13).Value = -71430929
$ws.Cells.Item(82, 14).Value = -3327.3333
$ws.Cells.Item(85, 8).Value = 38464204
$ws.Cells.Item(85, 9).Value = 71431290
$ws.Cells.Item(85, 10).Value = 2605.3333
$ws.Cells.Item(85, 11).Value = 71431290
$ws.Cells.Item(85, 12).Value = 2605.3333
$ws.Cells.Item(85, 13).Value = -71430042
$ws.Cells.Item(85, 14).Value = -5101.3333
$ws.Cells.Item(122, 8).Value = 5772.225
$ws.Cells.Item(122, 9).Value = 4558.3335
$ws.Cells.Item(122, 10).Value = 5986.4414
$ws.Cells.Item(122, 11).Value = 13675.0005
$ws.Cells.Item(122, 12).Value = 17959.3242
$ws.Cells.Item(122, 13).Value = -11225.0005
$ws.Cells.Item(122, 14).Value = -22859.3242
$ws.Cells.Item(124, 8).Value = 64214.5
$ws.Cells.Item(124, 10).Value = 64214.5
$ws.Cells.Item(124, 12).Value = 64214.5
$ws.Cells.Item(124, 14).Value = -74034.5
$ws.Cells.Item(126, 8).Value = 9741.286
$ws.Cells.Item(126, 9).Value = 9246
$ws.Cells.Item(126, 11).Value = 27738
$ws.Cells.Item(126, 13).Value = -25268
$ws.Cells.Item(132, 8).Value = 3312.05
$ws.Cells.Item(132, 9).Value = 2100.3
$ws.Cells.Item(132, 10).Value = 4523.8
$ws.Cells.Item(132, 11).Value = 6300.900000000001
$ws.Cells.Item(132, 12).Value = 13571.4
$ws.Cells.Item(132, 13).Value = -3770.900000000001
$ws.Cells.Item(132, 14).Value = -18631.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 947.6
$ws.Cells.Item(113, 9).Value = 1109.5
$ws.Cells.Item(113, 11).Value = 3328.5
$ws.Cells.Item(113, 13).Value = -1158.5
$ws.Cells.Item(126, 8).Value = 1595.4762
$ws.Cells.Item(126, 9).Value = 1713.7333
$ws.Cells.Item(126, 10).Value = 1299.8334
$ws.Cells.Item(126, 11).Value = 5141.199900000001
$ws.Cells.Item(126, 12).Value = 3899.5002
$ws.Cells.Item(126, 13).Value = -2671.199900000001
$ws.Cells.Item(126, 14).Value = -8839.5002
$ws.Cells.Item(132, 8).Value = 8335696.5
$ws.Cells.Item(132, 9).Value = 2248.5925
$ws.Cells.Item(132, 10).Value = 36461084
$ws.Cells.Item(132, 11).Value = 6745.7775
$ws.Cells.Item(132, 12).Value = 109383252
$ws.Cells.Item(132, 13).Value = -4215.7775
$ws.Cells.Item(132, 14).Value = -109388312
$ws.Cells.Item(136, 8).Value = 5201.2666
$ws.Cells.Item(136, 9).Value = 6284.5835
$ws.Cells.Item(136, 10).Value = 4479.0557
$ws.Cells.Item(136, 11).Value = 18853.7505
$ws.Cells.Item(136, 12).Value = 13437.1671
$ws.Cells.Item(136, 13).Value = -16303.7505
$ws.Cells.Item(136, 14).Value = -18537.1671
